$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.812.42'
$ws.Range("E2").Value = '  -0.45%  '

# Row 3
$ws.Range("D3").Value = '1.863.70'
$ws.Range("E3").Value = '  -1.55%  '

# Row 4
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.21%  '

# Row 5
$ws.Range("D5").Value = '''0.7346'
$ws.Range("E5").Value = '  -5.05%  '

# Row 6
$ws.Range("D6").Value = '''241.63'
$ws.Range("E6").Value = '  -0.84%  '

# Row 7
$ws.Range("D7").Value = '''1.003'
$ws.Range("E7").Value = '  +0.28%  '

# Row 8
$ws.Range("D8").Value = '''0.3087'
$ws.Range("E8").Value = '  -1.61%  '

# Row 9
$ws.Range("D9").Value = '''24.55'
$ws.Range("E9").Value = '  -4.64%  '

# Row 10
$ws.Range("D10").Value = '''0.07023'
$ws.Range("E10").Value = '  -4.53%  '

# Row 11
$ws.Range("D11").Value = '''0.08444'
$ws.Range("E11").Value = '  +4.69%  '

# Row 12
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '''0.7484'
$ws.Range("E12").Value = '  -2.99%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''5.305'
$ws.Range("E13").Value = '  -3.62%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.851.85'
$ws.Range("E14").Value = '  -1.66%  '

# Row 15
$ws.Range("D15").Value = '''92.04'
$ws.Range("E15").Value = '  -2.31%  '

# Row 16
$ws.Range("D16").Value = '29.784.48'
$ws.Range("E16").Value = '  -0.46%  '

# Row 17
$ws.Range("D17").Value = '''6.063'
$ws.Range("E17").Value = '  -2.54%  '

# Row 18
$ws.Range("D18").Value = '''13.50'
$ws.Range("E18").Value = '  -3.64%  '

# Row 19
$ws.Range("D19").Value = '''240.14'
$ws.Range("E19").Value = '  -2.68%  '

# Row 20
$ws.Range("D20").Value = '''0.000007757'
$ws.Range("E20").Value = '  -1.34%  '

# Row 21
$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.18%  '

# Row 22
$ws.Range("D22").Value = '2.140.28'
$ws.Range("E22").Value = '  +0.90%  '

# Row 23
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '''1.002'
$ws.Range("E23").Value = '  +0.16%  '

# Row 24
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '''7.895'
$ws.Range("E24").Value = '  -3.11%  '

# Row 25
$ws.Range("D25").Value = '''0.1561'
$ws.Range("E25").Value = '  -1.20%  '

# Row 26
$ws.Range("D26").Value = '''9.261'
$ws.Range("E26").Value = '  -1.92%  '

# Row 27
$ws.Range("D27").Value = '''162.62'
$ws.Range("E27").Value = '  +0.39%  '

# Row 28
$ws.Range("D28").Value = '''18.50'
$ws.Range("E28").Value = '  -1.43%  '

# Row 29
$ws.Range("D29").Value = '''2.001'
$ws.Range("E29").Value = '  -1.22%  '

# Row 30
$ws.Range("D30").Value = '''1.481'
$ws.Range("E30").Value = '  +4.06%  '

# Row 31
$ws.Range("D31").Value = '''1.528'
$ws.Range("E31").Value = '  -0.84%  '

# Row 32
$ws.Range("D32").Value = '''4.439'
$ws.Range("E32").Value = '  -0.66%  '

# Row 33
$ws.Range("D33").Value = '''4.125'
$ws.Range("E33").Value = '  +1.53%  '

# Row 34
$ws.Range("D34").Value = '''0.05352'
$ws.Range("E34").Value = '  -3.87%  '

# Row 35
$ws.Range("D35").Value = '''1.226'
$ws.Range("E35").Value = '  -1.10%  '

# Row 36
$ws.Range("D36").Value = '''0.7425'
$ws.Range("E36").Value = '  -1.28%  '

# Row 37
$ws.Range("D37").Value = '''1.005'
$ws.Range("E37").Value = '  +0.38%  '

# Row 38
$ws.Range("D38").Value = '''2.696'

# Row 39
$ws.Range("D39").Value = '''0.01924'
$ws.Range("E39").Value = '  -0.33%  '

# Row 40
$ws.Range("D40").Value = '''2.733'
$ws.Range("E40").Value = '  -2.04%  '

# Row 41
$ws.Range("D41").Value = '''0.4415'
$ws.Range("E41").Value = '  -1.31%  '

# Row 42
$ws.Range("D42").Value = '1.096.44'
$ws.Range("E42").Value = '  -0.48%  '

# Row 43
$ws.Range("D43").Value = '''6.020'
$ws.Range("E43").Value = '  +0.11%  '

# Row 44
$ws.Range("D44").Value = '''71.63'
$ws.Range("E44").Value = '  -3.81%  '

# Row 45
$ws.Range("D45").Value = '''0.8608'
$ws.Range("E45").Value = '  +1.12%  '

# Row 46
$ws.Range("D46").Value = '''1.004'
$ws.Range("E46").Value = '  +0.36%  '

# Row 47
$ws.Range("D47").Value = '''101.87'
$ws.Range("E47").Value = '  -0.52%  '

# Row 48
$ws.Range("D48").Value = '''7.643'
$ws.Range("E48").Value = '  +1.36%  '

# Row 49
$ws.Range("E49").Value = '  -3.29%  '

# Row 50
$ws.Range("D50").Value = '''3.003'
$ws.Range("E50").Value = '  +0.09%  '

# Row 51
$ws.Range("D51").Value = '2.037.88'
$ws.Range("E51").Value = '  -0.08%  '
